$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 33 corresponds to contest 21: PBKS vs KKR
# Fill in the "input" score columns (E, H, K, N, Q, T, W, Z, AC); the
# adjacent D/G/J/M/P/S/V/Y/AB columns already contain lookup formulas
# that will recalculate automatically.
$ws.Range("E33").Value = 70
$ws.Range("H33").Value = 80
$ws.Range("K33").Value = 100
$ws.Range("N33").Value = 50
$ws.Range("Q33").Value = 40
$ws.Range("T33").Value = 30
$ws.Range("W33").Value = 0
$ws.Range("Z33").Value = 60
$ws.Range("AC33").Value = 20

$excel.CalculateFull()
